$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- CasesTab secondary query (B2): append an ORDER BY / LIMIT clause ---
$b2 = $ws.Range("B2").Value2
$ws.Range("B2").Value = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100 "

# --- SamplesTab query (B3): append an ORDER BY / LIMIT clause ---
$b3 = $ws.Range("B3").Value2
$ws.Range("B3").Value = $b3 + "`n order By samp.sample_id ASC LIMIT 100"

# --- FilesTab query (B4): replace the trailing "order by f.file_name" with
#     "order By f.file_name ASC LIMIT 100" (note the extra leading space and
#     capitalised "By") ---
$b4 = $ws.Range("B4").Value2
$b4 = $b4 -replace "    order by f\.file_name$", "     order By f.file_name ASC LIMIT 100"
$ws.Range("B4").Value = $b4

# --- Row heights for rows 2 & 3 grew (Excel auto-fit after the text edit) ---
$ws.Rows.Item(2).RowHeight = 360
$ws.Rows.Item(3).RowHeight = 360

# --- Sheet view: selection moved from D4 to B4, and the view scrolled down
#     so row 3 is the top visible row ---
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 3
